$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 102400.52
$ws.Range("I15").Value = 102400.52
$ws.Range("K15").Value = 307201.56
$ws.Range("M15").Value = -307032.56
$ws.Range("H19").Value = 728.63635
$ws.Range("I19").Value = 531.8333
$ws.Range("J19").Value = 802.4375
$ws.Range("K19").Value = 531.8333
$ws.Range("L19").Value = 802.4375
$ws.Range("M19").Value = -356.8333
$ws.Range("N19").Value = -1152.4375
$ws.Range("H51").Value = 2417.389
$ws.Range("I51").Value = 2034.2858
$ws.Range("J51").Value = 2661.182
$ws.Range("K51").Value = 2034.2858
$ws.Range("L51").Value = 2661.182
$ws.Range("M51").Value = -1550.2858
$ws.Range("N51").Value = -3629.182
$ws.Range("H55").Value = 187.08333
$ws.Range("J55").Value = 277.2
$ws.Range("L55").Value = 277.2
$ws.Range("N55").Value = -705.2
$ws.Range("H69").Value = 4446544.5
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 11113361
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 33340083
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -33341831
$ws.Range("H70").Value = 2548.7
$ws.Range("I70").Value = 2201
$ws.Range("J70").Value = 2780.5
$ws.Range("K70").Value = 6603
$ws.Range("L70").Value = 8341.5
$ws.Range("M70").Value = -6333
$ws.Range("N70").Value = -8881.5
$ws.Range("H72").Value = 4446544.5
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 11113361
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 100020249
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -100028985
$ws.Range("H73").Value = 2548.7
$ws.Range("I73").Value = 2201
$ws.Range("J73").Value = 2780.5
$ws.Range("K73").Value = 6603
$ws.Range("L73").Value = 8341.5
$ws.Range("M73").Value = -5667
$ws.Range("N73").Value = -10213.5
$ws.Range("H76").Value = 4632717.5
$ws.Range("J76").Value = 3452.5
$ws.Range("L76").Value = 3452.5
$ws.Range("N76").Value = -4082.5
$ws.Range("H79").Value = 4632717.5
$ws.Range("J79").Value = 3452.5
$ws.Range("L79").Value = 3452.5
$ws.Range("N79").Value = -5636.5
$ws.Range("H87").Value = 25468
$ws.Range("J87").Value = 25666.666
$ws.Range("L87").Value = 25666.666
$ws.Range("N87").Value = -28162.666
$ws.Range("H90").Value = 25468
$ws.Range("J90").Value = 25666.666
$ws.Range("L90").Value = 76999.99800000001
$ws.Range("N90").Value = -89479.99800000001
$ws.Range("H107").Value = 2222692.8
$ws.Range("I107").Value = 2222692.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2222692.8
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -2220772.8
$ws.Range("H123").Value = 97591.336
$ws.Range("J123").Value = 97591.336
$ws.Range("L123").Value = 97591.336
$ws.Range("N123").Value = -107391.336
$ws.Range("H131").Value = 4087.087
$ws.Range("I131").Value = 876.64703
$ws.Range("J131").Value = 13183.333
$ws.Range("K131").Value = 2629.94109
$ws.Range("L131").Value = 39549.999
$ws.Range("M131").Value = 2410.05891
$ws.Range("N131").Value = -49629.999
$ws.Range("H132").Value = 307180.3
$ws.Range("I132").Value = 419734.06
$ws.Range("J132").Value = 35175.418
$ws.Range("K132").Value = 1259202.18
$ws.Range("L132").Value = 105526.254
$ws.Range("M132").Value = -1256672.18
$ws.Range("N132").Value = -110586.254
$ws.Range("H137").Value = 1927.0588
$ws.Range("I137").Value = 1800
$ws.Range("K137").Value = 5400
$ws.Range("M137").Value = -2850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2621.3142
$ws.Range("I32").Value = 1862.7903
$ws.Range("K32").Value = 1862.7903
$ws.Range("M32").Value = -1575.7903
$ws.Range("H45").Value = 2505.2
$ws.Range("I45").Value = 1618
$ws.Range("J45").Value = 4279.6
$ws.Range("K45").Value = 1618
$ws.Range("L45").Value = 4279.6
$ws.Range("M45").Value = -1241
$ws.Range("N45").Value = -5033.6
$ws.Range("H97").Value = 18524526
$ws.Range("I97").Value = 19614162
$ws.Range("J97").Value = 711
$ws.Range("K97").Value = 19614162
$ws.Range("L97").Value = 711
$ws.Range("M97").Value = -19613666
$ws.Range("N97").Value = -1703
$ws.Range("H110").Value = 862.5294
$ws.Range("I110").Value = 790
$ws.Range("J110").Value = 1406.5
$ws.Range("K110").Value = 790
$ws.Range("L110").Value = 1406.5
$ws.Range("M110").Value = 1255
$ws.Range("N110").Value = -5496.5
$ws.Range("H122").Value = 1567.6
$ws.Range("J122").Value = 1669.909
$ws.Range("L122").Value = 5009.727000000001
$ws.Range("N122").Value = -9909.727000000001
$ws.Range("H123").Value = 33609.332
$ws.Range("J123").Value = 33609.332
$ws.Range("L123").Value = 33609.332
$ws.Range("N123").Value = -43409.332
$ws.Range("H131").Value = 44900
$ws.Range("J131").Value = 44900
$ws.Range("L131").Value = 44900
$ws.Range("N131").Value = -54980
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1414.5
$ws.Range("I94").Value = 1458
$ws.Range("J94").Value = 1110
$ws.Range("K94").Value = 1458
$ws.Range("L94").Value = 1110
$ws.Range("M94").Value = -1007
$ws.Range("N94").Value = -2012
$ws.Range("H99").Value = 3170.8333
$ws.Range("I99").Value = 1640.909
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 1640.909
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = -142.9090000000001
$ws.Range("N99").Value = -22996
$ws.Range("H105").Value = 12348383
$ws.Range("I105").Value = 13891514
$ws.Range("J105").Value = 3333.3333
$ws.Range("K105").Value = 13891514
$ws.Range("L105").Value = 3333.3333
$ws.Range("M105").Value = -13889767
$ws.Range("N105").Value = -6827.3333
$ws.Range("H107").Value = 1124.619
$ws.Range("I107").Value = 938.6111
$ws.Range("J107").Value = 2240.6667
$ws.Range("K107").Value = 938.6111
$ws.Range("L107").Value = 2240.6667
$ws.Range("M107").Value = 981.3889
$ws.Range("N107").Value = -6080.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1409
$ws.Range("I122").Value = 1118.4
$ws.Range("J122").Value = 1699.6
$ws.Range("K122").Value = 3355.2
$ws.Range("L122").Value = 5098.799999999999
$ws.Range("M122").Value = -905.2000000000003
$ws.Range("N122").Value = -9998.799999999999
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("N131").Value = 0
$ws.Range("H132").Value = 3724.2
$ws.Range("I132").Value = 1325.5
$ws.Range("K132").Value = 3976.5
$ws.Range("M132").Value = -1446.5
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1069.3673
$ws.Range("I68").Value = 853.84485
$ws.Range("J68").Value = 1381.875
$ws.Range("K68").Value = 2561.53455
$ws.Range("L68").Value = 4145.625
$ws.Range("M68").Value = -1750.53455
$ws.Range("N68").Value = -5767.625
$ws.Range("H71").Value = 1069.3673
$ws.Range("I71").Value = 853.84485
$ws.Range("J71").Value = 1381.875
$ws.Range("K71").Value = 7684.603649999999
$ws.Range("L71").Value = 12436.875
$ws.Range("M71").Value = -3628.603649999999
$ws.Range("N71").Value = -20548.875
$ws.Range("H127").Value = 3166.8667
$ws.Range("J127").Value = 3166.8667
$ws.Range("L127").Value = 9500.6001
$ws.Range("N127").Value = -19420.6001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6405.6313
$ws.Range("I70").Value = 6940.467
$ws.Range("J70").Value = 4400
$ws.Range("K70").Value = 6940.467
$ws.Range("L70").Value = 4400
$ws.Range("M70").Value = -6670.467
$ws.Range("N70").Value = -4940
$ws.Range("H73").Value = 6405.6313
$ws.Range("I73").Value = 6940.467
$ws.Range("J73").Value = 4400
$ws.Range("K73").Value = 6940.467
$ws.Range("L73").Value = 4400
$ws.Range("M73").Value = -6004.467
$ws.Range("N73").Value = -6272
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("N91").Value = 0
$ws.Range("H102").Value = 4476.222
$ws.Range("I102").Value = 1680.4
$ws.Range("J102").Value = 7971
$ws.Range("K102").Value = 1680.4
$ws.Range("L102").Value = 7971
$ws.Range("M102").Value = -58.40000000000009
$ws.Range("N102").Value = -11215
$ws.Range("H132").Value = 2484.1462
$ws.Range("I132").Value = 1851.5151
$ws.Range("J132").Value = 5093.75
$ws.Range("K132").Value = 5554.5453
$ws.Range("L132").Value = 15281.25
$ws.Range("M132").Value = -3024.5453
$ws.Range("N132").Value = -20341.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3992.5
$ws.Range("I122").Value = 3990
$ws.Range("J122").Value = 3992.8572
$ws.Range("K122").Value = 11970
$ws.Range("L122").Value = 11978.5716
$ws.Range("M122").Value = -9520
$ws.Range("N122").Value = -16878.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 201380.8
$ws.Range("I122").Value = 501002
$ws.Range("K122").Value = 1503006
$ws.Range("M122").Value = -1500556
$ws.Range("H123").Value = 33283.668
$ws.Range("J123").Value = 33283.668
$ws.Range("L123").Value = 33283.668
$ws.Range("N123").Value = -43083.668
$ws.Range("H132").Value = 23812822
$ws.Range("I132").Value = 31252776
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 93758328
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -93755798
$ws.Range("N132").Value = -19964
$ws.Range("H136").Value = 19668284
$ws.Range("I136").Value = 25718568
$ws.Range("J136").Value = 4863.25
$ws.Range("K136").Value = 77155704
$ws.Range("L136").Value = 14589.75
$ws.Range("M136").Value = -77153154
$ws.Range("N136").Value = -19689.75
